$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 7000
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -7968

$ws.Range("H55").Value = 462.33334
$ws.Range("I55").Value = 309
$ws.Range("J55").Value = 539
$ws.Range("K55").Value = 309
$ws.Range("L55").Value = 539
$ws.Range("M55").Value = -95
$ws.Range("N55").Value = -967

$ws.Range("H88").Value = 4250
$ws.Range("J88").Value = 4250
$ws.Range("L88").Value = 4250
$ws.Range("N88").Value = -5062

$ws.Range("H91").Value = 4250
$ws.Range("J91").Value = 4250
$ws.Range("L91").Value = 4250
$ws.Range("N91").Value = -7058

$ws.Range("H106").Value = 1499
$ws.Range("I106").Value = 1499
$ws.Range("K106").Value = 1499
$ws.Range("M106").Value = -868

$ws.Range("H137").Value = 1567.3214
$ws.Range("I137").Value = 1314.8462
$ws.Range("J137").Value = 4849.5
$ws.Range("K137").Value = 3944.5386
$ws.Range("L137").Value = 14548.5
$ws.Range("M137").Value = -1394.5386
$ws.Range("N137").Value = -19648.5

$ws.Range("H138").Value = 3410.4
$ws.Range("J138").Value = 5199.6665
$ws.Range("L138").Value = 15598.9995
$ws.Range("N138").Value = -25878.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 415.8
$ws.Range("I26").Value = 501
$ws.Range("J26").Value = 75
$ws.Range("K26").Value = 501
$ws.Range("L26").Value = 75
$ws.Range("M26").Value = -171
$ws.Range("N26").Value = -735

$ws.Range("H74").Value = 1004.5
$ws.Range("I74").Value = 1004.5
$ws.Range("K74").Value = 1004.5
$ws.Range("M74").Value = -130.5

$ws.Range("H77").Value = 1004.5
$ws.Range("I77").Value = 1004.5
$ws.Range("K77").Value = 5022.5
$ws.Range("M77").Value = -654.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5414.4
$ws.Range("I86").Value = 4681.6
$ws.Range("K86").Value = 4681.6
$ws.Range("M86").Value = -3558.6

$ws.Range("H89").Value = 5414.4
$ws.Range("I89").Value = 4681.6
$ws.Range("K89").Value = 23408
$ws.Range("M89").Value = -17792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2796.3
$ws.Range("J31").Value = 2517.8
$ws.Range("L31").Value = 2517.8
$ws.Range("N31").Value = -3107.8

$ws.Range("H32").Value = 992.25
$ws.Range("I32").Value = 489.66666
$ws.Range("K32").Value = 489.66666
$ws.Range("M32").Value = -173.66666

$ws.Range("H34").Value = 2796.3
$ws.Range("J34").Value = 2517.8
$ws.Range("L34").Value = 2517.8
$ws.Range("N34").Value = -2921.8

$ws.Range("H132").Value = 2041.12
$ws.Range("I132").Value = 2083.8262
$ws.Range("K132").Value = 6251.4786
$ws.Range("M132").Value = -3721.4786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 25249.5
$ws.Range("J47").Value = 25249.5
$ws.Range("L47").Value = 25249.5
$ws.Range("N47").Value = -26385.5

$ws.Range("H48").Value = 37000
$ws.Range("J48").Value = 37000
$ws.Range("L48").Value = 37000
$ws.Range("N48").Value = -37970

$ws.Range("H55").Value = 9997
$ws.Range("I55").Value = 9997
$ws.Range("K55").Value = 9997
$ws.Range("M55").Value = -9670

$ws.Range("H80").Value = 2979.5
$ws.Range("J80").Value = 3500
$ws.Range("L80").Value = 3500
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 2979.5
$ws.Range("J83").Value = 3500
$ws.Range("L83").Value = 17500
$ws.Range("N83").Value = -27484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2214.2856
$ws.Range("I22").Value = 1375
$ws.Range("J22").Value = 3333.3333
$ws.Range("K22").Value = 1375
$ws.Range("L22").Value = 3333.3333
$ws.Range("M22").Value = -1080
$ws.Range("N22").Value = -3923.3333

$ws.Range("H27").Value = 2214.2856
$ws.Range("I27").Value = 1375
$ws.Range("J27").Value = 3333.3333
$ws.Range("K27").Value = 1375
$ws.Range("L27").Value = 3333.3333
$ws.Range("M27").Value = -1268
$ws.Range("N27").Value = -3547.3333

$ws.Range("H42").Value = 35000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 35000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 35000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -36126

$ws.Range("H43").Value = 24999.5
$ws.Range("J43").Value = 24999.5
$ws.Range("L43").Value = 24999.5
$ws.Range("N43").Value = -25385.5

$ws.Range("H49").Value = 35000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 35000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 35000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -35294

$ws.Range("H68").Value = 3751.5
$ws.Range("I68").Value = 3500
$ws.Range("K68").Value = 3500
$ws.Range("M68").Value = -2751

$ws.Range("H69").Value = 55555
$ws.Range("J69").Value = 55555
$ws.Range("L69").Value = 55555
$ws.Range("N69").Value = -57177

$ws.Range("H71").Value = 3751.5
$ws.Range("I71").Value = 3500
$ws.Range("K71").Value = 17500
$ws.Range("M71").Value = -13756

$ws.Range("H72").Value = 55555
$ws.Range("J72").Value = 55555
$ws.Range("L72").Value = 166665
$ws.Range("N72").Value = -174777

$ws.Range("H140").Value = 75833.336
$ws.Range("I140").Value = 72500
$ws.Range("J140").Value = 77500
$ws.Range("K140").Value = 72500
$ws.Range("L140").Value = 77500
$ws.Range("M140").Value = -67320
$ws.Range("N140").Value = -87860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19480
